$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bibi Cell Mundi)
$ws.Range("AC2").Value = 13062.18
$ws.Range("AD2").Value = 12545.8
$ws.Range("AG2").Value = 261462.34

# Row 3 (Bibi Cell Vieiralves)
$ws.Range("AC3").Value = 4963
$ws.Range("AD3").Value = 4496
$ws.Range("AG3").Value = 106704.05

# Row 4 (Bibi Cell Manauara)
$ws.Range("AC4").Value = 3823
$ws.Range("AD4").Value = 2304
$ws.Range("AG4").Value = 94604.5

# Row 5 (Bibi Cell Ponta Negra)
$ws.Range("AC5").Value = 3832
$ws.Range("AD5").Value = 2631
$ws.Range("AG5").Value = 79901.84

# Row 6 (total)
$ws.Range("AC6").Value = 25680.18
$ws.Range("AD6").Value = 21976.8
$ws.Range("AG6").Value = 542672.73
